$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 19 (old rows 23,24 shift to 27,28)
$ws.Rows("19:22").Insert()

# The old row 18 carried the table's "bottom border" styling (it was the
# last data row). After the insert it becomes a normal middle row, and the
# new row 22 (the new last data row) should carry that bottom-border style
# instead. Re-stripe: rows 18-21 get the "middle" style (copied from row 17),
# and row 22 gets the old "bottom" style (copied from the original row 18,
# which - pre-insert - is still intact at this point... but since Insert
# already happened, grab the bottom style from row 18 BEFORE overwriting it).
$ws.Range("B18:J18").Copy() | Out-Null
$ws.Range("B22:J22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B17:J17").Copy() | Out-Null
$ws.Range("B18:J21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Fill rows 16-18 with period 2508 (unchanged data), rows 19-22 with period 2509
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143383078"
$ws.Range("D16").Value = "LAURA ISABEL PUERTA PIZARRO"
$ws.Range("E16").Value = "2508"
$ws.Range("F16").Value = 106160
$ws.Range("G16").Value = 2654000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1050968160"
$ws.Range("D17").Value = "RAUL ANDRES BORRERO BERMUDEZ"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 75360
$ws.Range("G17").Value = 1884000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1007314724"
$ws.Range("D18").Value = "MIGUEL JOSE GOMEZ PEREZ"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 75360
$ws.Range("G18").Value = 1884000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143383078"
$ws.Range("D19").Value = "LAURA ISABEL PUERTA PIZARRO"
$ws.Range("E19").Value = "2509"
$ws.Range("F19").Value = 106160
$ws.Range("G19").Value = 2654000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1050968160"
$ws.Range("D20").Value = "RAUL ANDRES BORRERO BERMUDEZ"
$ws.Range("E20").Value = "2509"
$ws.Range("F20").Value = 75360
$ws.Range("G20").Value = 1884000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1007314724"
$ws.Range("D21").Value = "MIGUEL JOSE GOMEZ PEREZ"
$ws.Range("E21").Value = "2509"
$ws.Range("F21").Value = 75360
$ws.Range("G21").Value = 1884000

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1002431293"
$ws.Range("D22").Value = "LUIS DANIEL SANDOVAL MUÑIZ"
$ws.Range("E22").Value = "2509"
$ws.Range("F22").Value = 59840
$ws.Range("G22").Value = 1496000

# Update summary fields
$ws.Range("E11").Value = 573600
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 2
